# "Generate Report for Handoff"
# Updates the localization-status report: flips the "Handed back: in sync
# with en-US" status to "Ready for handoff" everywhere it appears, bumps
# the handoff-generation timestamps to reflect the new run, and re-fits
# the now-shorter "Status" columns.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldDate1 = "2016-09-04 11:04:25"
$newDate1 = "2016-09-04 11:05:13"

$oldDate2 = "2016-09-04 11:04:21"
$newDate2 = "2016-09-04 11:05:05"

# New, tighter width for the Status columns now that the text is shorter.
$newStatusColWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newDate1

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newDate2

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newDate1

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
